# Updated cryptos list with GitHub Actions (refresh of Price / Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = "67.157.49"
$ws.Cells.Item(2, 5).Value = "  -1.50%  "

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = "3.487.56"
$ws.Cells.Item(3, 5).Value = "  -4.01%  "

# Row 4: TetherUSD
$ws.Cells.Item(4, 5).Value = "  -0.10%  "

# Row 5: Solana
$ws.Cells.Item(5, 4).Value = "'198.96"
$ws.Cells.Item(5, 5).Value = "  +1.64%  "

# Row 6: BNB
$ws.Cells.Item(6, 4).Value = "'548.16"
$ws.Cells.Item(6, 5).Value = "  -4.52%  "

# Row 7: LidoStakedEther
$ws.Cells.Item(7, 4).Value = "3.484.80"
$ws.Cells.Item(7, 5).Value = "  -3.96%  "

# Row 8: XRP
$ws.Cells.Item(8, 4).Value = "'0.604"
$ws.Cells.Item(8, 5).Value = "  -2.68%  "

# Row 9: USDC
$ws.Cells.Item(9, 5).Value = "  -0.08%  "

# Row 10: Cardano
$ws.Cells.Item(10, 4).Value = "'0.649"
$ws.Cells.Item(10, 5).Value = "  -4.25%  "

# Row 11: Avalanche
$ws.Cells.Item(11, 4).Value = "'61.96"
$ws.Cells.Item(11, 5).Value = "  +10.36%  "

# Row 12: Dogecoin
$ws.Cells.Item(12, 4).Value = "'0.142"
$ws.Cells.Item(12, 5).Value = "  -7.03%  "

# Row 13: ShibaInu
$ws.Cells.Item(13, 4).Value = "'0.0000267"
$ws.Cells.Item(13, 5).Value = "  -9.03%  "

# Row 14: Polkadot
$ws.Cells.Item(14, 4).Value = "'9.72"
$ws.Cells.Item(14, 5).Value = "  -3.96%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Cells.Item(15, 4).Value = "4.032.91"
$ws.Cells.Item(15, 5).Value = "  -4.17%  "

# Row 16: WrappedEther
$ws.Cells.Item(16, 4).Value = "3.479.76"
$ws.Cells.Item(16, 5).Value = "  -4.42%  "

# Row 17: TRON
$ws.Cells.Item(17, 5).Value = "  -2.11%  "

# Row 18: WrappedBTC
$ws.Cells.Item(18, 4).Value = "66.775.97"
$ws.Cells.Item(18, 5).Value = "  -1.96%  "

# Row 19: Chainlink
$ws.Cells.Item(19, 4).Value = "'18.15"
$ws.Cells.Item(19, 5).Value = "  -2.10%  "

# Row 20: Uniswap
$ws.Cells.Item(20, 4).Value = "'11.69"
$ws.Cells.Item(20, 5).Value = "  -6.84%  "

# Row 21: Polygon
$ws.Cells.Item(21, 4).Value = "'1.02"
$ws.Cells.Item(21, 5).Value = "  -5.77%  "

# Row 22: BitcoinCash
$ws.Cells.Item(22, 4).Value = "'387.08"
$ws.Cells.Item(22, 5).Value = "  -3.85%  "

# Row 23: PancakeSwap
$ws.Cells.Item(23, 2).Value = "RenderToken"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(23, 4).Value = "'11.91"
$ws.Cells.Item(23, 5).Value = "  -6.42%  "

# Row 24: RenderToken
$ws.Cells.Item(24, 2).Value = "PancakeSwap"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(24, 4).Value = "'3.97"
$ws.Cells.Item(24, 5).Value = "  -5.96%  "

# Row 25: Litecoin
$ws.Cells.Item(25, 4).Value = "'82.08"
$ws.Cells.Item(25, 5).Value = "  -4.57%  "

# Row 26: Toncoin
$ws.Cells.Item(26, 4).Value = "'3.83"
$ws.Cells.Item(26, 5).Value = "  -0.59%  "

# Row 27: InternetComputer(DFINITY)
$ws.Cells.Item(27, 4).Value = "'12.09"
$ws.Cells.Item(27, 5).Value = "  -4.12%  "

# Row 28: ImmutableX
$ws.Cells.Item(28, 4).Value = "'2.77"
$ws.Cells.Item(28, 5).Value = "  -5.99%  "

# Row 29: Filecoin
$ws.Cells.Item(29, 4).Value = "'8.75"
$ws.Cells.Item(29, 5).Value = "  -4.39%  "

# Row 30: EthereumClassic
$ws.Cells.Item(30, 4).Value = "'30.84"
$ws.Cells.Item(30, 5).Value = "  -2.62%  "

# Row 31: Bittensor
$ws.Cells.Item(31, 4).Value = "'673.67"
$ws.Cells.Item(31, 5).Value = "  -2.81%  "

# Row 32: NEARProtocol
$ws.Cells.Item(32, 4).Value = "'6.91"
$ws.Cells.Item(32, 5).Value = "  -14.78%  "

# Row 33: Cosmos
$ws.Cells.Item(33, 4).Value = "'11.62"
$ws.Cells.Item(33, 5).Value = "  -4.81%  "

# Row 34: OKB
$ws.Cells.Item(34, 4).Value = "'63.54"
$ws.Cells.Item(34, 5).Value = "  -1.96%  "

# Row 35: Hedera
$ws.Cells.Item(35, 5).Value = "  -7.61%  "

# Row 36: InjectiveProtocol
$ws.Cells.Item(36, 4).Value = "'38.32"
$ws.Cells.Item(36, 5).Value = "  -10.24%  "

# Row 37: Dai
$ws.Cells.Item(37, 5).Value = "  -0.01%  "

# Row 38: TheGraph
$ws.Cells.Item(38, 4).Value = "'0.394"
$ws.Cells.Item(38, 5).Value = "  -5.28%  "

# Row 39: FirstDigitalUSD
$ws.Cells.Item(39, 4).Value = "'0.997"
$ws.Cells.Item(39, 5).Value = "  -0.22%  "

# Row 40: Maker
$ws.Cells.Item(40, 4).Value = "3.055.79"
$ws.Cells.Item(40, 5).Value = "  -3.86%  "

# Row 41: Kaspa
$ws.Cells.Item(41, 4).Value = "'0.129"
$ws.Cells.Item(41, 5).Value = "  -5.03%  "

# Row 42: ThetaToken
$ws.Cells.Item(42, 4).Value = "'2.96"
$ws.Cells.Item(42, 5).Value = "  -5.00%  "

# Row 43: PEPE
$ws.Cells.Item(43, 4).Value = "0.0₃0669"
$ws.Cells.Item(43, 5).Value = "  -16.04%  "

# Row 44: WEMIXToken
$ws.Cells.Item(44, 4).Value = "'2.76"
$ws.Cells.Item(44, 5).Value = "  +6.25%  "

# Row 45: Fetch.AI
$ws.Cells.Item(45, 4).Value = "'2.48"
$ws.Cells.Item(45, 5).Value = "  -13.21%  "

# Row 46: dogwifhat
$ws.Cells.Item(46, 4).Value = "'2.71"
$ws.Cells.Item(46, 5).Value = "  -6.90%  "

# Row 47: VeChain
$ws.Cells.Item(47, 4).Value = "'0.0394"
$ws.Cells.Item(47, 5).Value = "  -6.60%  "

# Row 48: Stellar
$ws.Cells.Item(48, 4).Value = "'0.126"
$ws.Cells.Item(48, 5).Value = "  -4.97%  "

# Row 49: Monero
$ws.Cells.Item(49, 4).Value = "'135.85"
$ws.Cells.Item(49, 5).Value = "  -5.01%  "

# Row 50: THORChain
$ws.Cells.Item(50, 4).Value = "'8.13"
$ws.Cells.Item(50, 5).Value = "  -8.18%  "

# Row 51: ApeXProtocol
$ws.Cells.Item(51, 4).Value = "'2.85"
$ws.Cells.Item(51, 5).Value = "  -8.59%  "
